# Apply crypto price/volume updates from the commit diff.
# Pure-numeric-looking strings in column D must stay TEXT (not be
# coerced to floating point numbers by Excel), so those cells get a
# temporary Text number format around the assignment, then are reset
# back to the default "Normal" style so no stray formatting remains.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.080.24"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.820.93"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.81"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.20%  "
$ws.Range("E9").Value = "  +4.66%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0991"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "2.083.93"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "1.817.90"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "35.062.33"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.124"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.19%  "
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.698"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "92.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.40%  "
$ws.Range("D39").Value = "1.341.81"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.987"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").Value = "2.000.60"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +4.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.53%  "
